# Remove PCA9306 level shifter (U6 U7) row from the BOM sheet.
# The 9306 level shifter isn't used in the design, so it shouldn't be in the BOM.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOM")

# Row 35 contains the "U6 U7" / "PCA9306" bill-of-materials line.
# Deleting the entire row shifts every following row up by one.
$ws.Rows.Item(35).Delete()
